$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style used by
# the existing header row (e.g. H1).
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill data rows 2-33: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($r, 8).Value2
}

$ws.Range("A1").Select()
